$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 170234
$ws.Range("C4").Value = 161057
$ws.Range("C5").Value = 9177
$ws.Range("C8").Value = 65.79000000000001
